$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary Sheet")

# Expand the print area from $B$1:$J$17 to $B$1:$K$17
$ws.PageSetup.PrintArea = "`$B`$1:`$K`$17"

# Update the selection to B1:K17 (was B1:K1)
$ws.Range("B1:K17").Select()

# Reduce the print scale from 78% to 70%
$ws.PageSetup.Zoom = 70
